$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each written cell to remain Text (matches the source inlineStr cells),
# then reset the style back to Normal so no stray number-format/quote-prefix
# style survives the write (Excel otherwise auto-detects "6.20" etc. as numbers).

$ws.Range('D2').Value = "'43.641.89"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.27%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.287.93"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.21%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.15%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'95.25"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -4.18%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'267.34"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.65%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.622"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -1.00%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.00%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.608"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -4.11%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'44.61"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -7.84%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0937"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -1.12%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'7.76"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -5.51%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.58%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'2.631.85"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.29%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'15.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -2.61%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.847"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.83%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'2.289.18"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +1.55%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'43.580.10"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.29%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.37%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'6.20"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.77%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'72.37"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.78%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'2.45"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +3.94%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'234.82"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.14%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'9.02"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -15.50%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.05%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.50"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.27%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'11.21"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.67%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'3.46"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.07%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'40.46"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.74%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.13%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'175.35"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.02%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'21.90"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +3.12%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.0881"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -4.33%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'5.34"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -6.82%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.22%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value = "'VeChain"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'0.0356"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.17%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = "'Kaspa"
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'0.108"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -6.13%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'4.39"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.62%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -7.64%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +6.71%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  -7.49%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +15.37%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'64.42"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +2.63%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'12.04"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -4.85%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'8.82"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +2.62%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'5.23"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -4.87%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.25%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'98.01"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.62%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.04%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.512.05"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.38%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +4.59%  "
$ws.Range('E51').Style = 'Normal'
